$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "random_key"
$ws.Range("B5").Value = "abc"
$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "string"

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Adding the new row changes the natural wrapped height of rows 2 & 4,
# so Excel re-autofits them (dropping their stale explicit row heights).
$ws.Range("A2:D2").EntireRow.AutoFit()
$ws.Range("A4:D4").EntireRow.AutoFit()

$ws.Range("C7").Select() | Out-Null
